# The "reviews_count" column (E) is empty for every row and is being
# removed entirely, shifting the subsequent columns (reviews_average,
# latitude, longitude, is_permanently_closed, gmaps_link,
# latest_review_date) one place to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E:E").Delete()
